# Avance 1 mod 1
# Modificación en avance Sebastián , creación de texto1
#
# Adds a new bulleted "Prrafodelista" (List Paragraph) item, right after
# the existing "Propuesta de Arquitectura de SW..." bullet and before the
# trailing blank paragraph, containing the text "dasdasdasdas".

$d = $word.ActiveDocument

# Locate the last bullet item of the numbered list (numId 4), i.e. the
# paragraph that ends with "...desventajas)." so the new paragraph can be
# inserted right after it, inheriting the same paragraph/run formatting.
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*desventajas)*") {
        $anchorPara = $p
    }
}

# Insert a brand-new paragraph right after the anchor paragraph; Word
# automatically carries over the paragraph style / numbering / run
# formatting from the paragraph it follows.
$anchorPara.Range.InsertParagraphAfter()

$newParaIndex = $anchorPara.Index + 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newPara.Range.Text = "dasdasdasdas"
